$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 265; existing rows 265-310 shift down to 266-311,
# and the used range grows from A1:R310 to A1:R311.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new daily price record.
$ws.Range("A265").Value = 5
$ws.Range("B265").Value = "Macroferia Regional de Talca"
$ws.Range("C265").Value = "Maule"
$ws.Range("D265").Value = "2021-11-04"
$ws.Range("E265").Value = 7
$ws.Range("F265").Value = 100114001
$ws.Range("G265").Value = "Papa"
$ws.Range("H265").Value = "Rodeo"
$ws.Range("I265").Value = "1a nueva(o)"
$ws.Range("J265").Value = 1500
$ws.Range("K265").Value = 10000
$ws.Range("L265").Value = 10000
$ws.Range("M265").Value = 10000
$ws.Range("N265").Value = "$/saco 25 kilos"
$ws.Range("O265").Value = "Región de O'Higgins"
$ws.Range("P265").Value = 400
$ws.Range("Q265").Value = 25
$ws.Range("R265").Value = "Hortaliza"
